$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("B2").Formula = "=4700*2"
$ws.Range("D2").Value = 5
$ws.Range("F2").Formula = "=((B2/(A2+B2))*D2)-E2"
$ws.Range("N2").Formula = "=F2/2"

# New Row 4
$ws.Range("A4").Formula = "=47000*2"
$ws.Range("B4").Formula = "=10000*2"
$ws.Range("D4").Value = 19
$ws.Range("E4").Value = 0
$ws.Range("F4").Formula = "=((B4/(A4+B4))*D4)-E4"
$ws.Range("H4").Formula = "=D4/(A4+B4)"
$ws.Range("K4").Formula = "=A4+B4"
$ws.Range("L4").Formula = "=B4/K4"
$ws.Range("N4").Formula = "=F4/2"

# New Row 7
$ws.Range("H7").Formula = "=2*H2+H4"

# New Row 10
$ws.Range("D10").Formula = "=19/4096"

# Selection
$ws.Range("D24").Select() | Out-Null
